$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 146, pushing existing rows 146-215 down to 147-216
$ws.Rows.Item(146).Insert()

# Fill in the constant columns (A-C, E-L) identical to the row that used to be here
$ws.Cells.Item(146, 1).Value = 7
$ws.Cells.Item(146, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(146, 3).Value = "Ñuble"
$ws.Cells.Item(146, 4).Value = 44879
$ws.Cells.Item(146, 5).Value = 16
$ws.Cells.Item(146, 6).Value = "Fruta"
$ws.Cells.Item(146, 7).Value = 100102
$ws.Cells.Item(146, 8).Value = "Cítricos"
$ws.Cells.Item(146, 9).Value = 100102004
$ws.Cells.Item(146, 10).Value = "Mandarina"
$ws.Cells.Item(146, 11).Value = "Murcott"
$ws.Cells.Item(146, 12).Value = "Primera"
$ws.Cells.Item(146, 13).Value = 120
$ws.Cells.Item(146, 14).Value = 8000
$ws.Cells.Item(146, 15).Value = 8000
$ws.Cells.Item(146, 16).Value = 8000
$ws.Cells.Item(146, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(146, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(146, 19).Value = 533
$ws.Cells.Item(146, 20).Value = 15

# Match the date-formatted style used by column D in the rest of the table
$ws.Cells.Item(146, 4).NumberFormat = $ws.Cells.Item(147, 4).NumberFormat
